$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.654.87'
$ws.Range("E2").Value = '  -0.29%  '

$ws.Range("D3").Value = '1.689.08'
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("E4").Value = '  +1.04%  '

$ws.Range("D5").Value = "'316.06"
$ws.Range("E5").Value = '  +1.18%  '

$ws.Range("E6").Value = '  +1.03%  '

$ws.Range("D7").Value = "'0.3938"
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = "'0.4045"
$ws.Range("E8").Value = '  -0.38%  '

$ws.Range("E9").Value = '  -2.13%  '

$ws.Range("D10").Value = "'1.001"
$ws.Range("E10").Value = '  +0.99%  '

$ws.Range("D11").Value = "'52.91"
$ws.Range("E11").Value = '  -1.81%  '

$ws.Range("D12").Value = "'0.08804"

$ws.Range("D13").Value = "'7.218"
$ws.Range("E13").Value = '  -1.11%  '

$ws.Range("D14").Value = "'23.48"
$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").Value = "'8.067"
$ws.Range("E15").Value = '  +8.11%  '

$ws.Range("E16").Value = '  -0.69%  '

$ws.Range("D17").Value = '1.695.78'
$ws.Range("E17").Value = '  +0.09%  '

$ws.Range("D18").Value = "'99.66"
$ws.Range("E18").Value = '  -0.85%  '

$ws.Range("D19").Value = "'0.07001"
$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("D20").Value = "'19.45"
$ws.Range("E20").Value = '  -0.16%  '

$ws.Range("D21").Value = "'7.006"
$ws.Range("E21").Value = '  +3.87%  '

$ws.Range("E22").Value = '  +1.10%  '

$ws.Range("D23").Value = "'14.28"
$ws.Range("E23").Value = '  +0.79%  '

$ws.Range("D24").Value = '24.639.32'
$ws.Range("E24").Value = '  -0.36%  '

$ws.Range("D25").Value = "'3.273"
$ws.Range("E25").Value = '  +10.37%  '

$ws.Range("E26").Value = '  +2.67%  '

$ws.Range("E27").Value = '  +1.24%  '

$ws.Range("D28").Value = "'162.78"
$ws.Range("E28").Value = '  +2.55%  '

$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = "'135.21"
$ws.Range("E29").Value = '  +1.75%  '

$ws.Range("B30").Value = 'HuobiToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D30").Value = "'5.173"
$ws.Range("E30").Value = '  +1.04%  '

$ws.Range("D31").Value = "'7.626"
$ws.Range("E31").Value = '  +2.78%  '

$ws.Range("D32").Value = '1.881.06'
$ws.Range("E32").Value = '  +0.07%  '

$ws.Range("D33").Value = "'0.08547"
$ws.Range("E33").Value = '  -1.67%  '

$ws.Range("D34").Value = "'1.055"
$ws.Range("E34").Value = '  -4.27%  '

$ws.Range("D35").Value = "'7.193"
$ws.Range("E35").Value = '  -2.22%  '

$ws.Range("D36").Value = "'11.11"
$ws.Range("E36").Value = '  -0.26%  '

$ws.Range("D37").Value = "'0.2732"
$ws.Range("E37").Value = '  +0.10%  '

$ws.Range("D38").Value = "'1.896"
$ws.Range("E38").Value = '  -1.67%  '

$ws.Range("D39").Value = "'14.30"
$ws.Range("E39").Value = '  -3.51%  '

$ws.Range("D40").Value = "'0.09170"
$ws.Range("E40").Value = '  +2.40%  '

$ws.Range("E41").Value = '  -2.92%  '

$ws.Range("D42").Value = "'1.457"
$ws.Range("E42").Value = '  -1.22%  '

$ws.Range("D43").Value = "'0.7588"
$ws.Range("E43").Value = '  -0.67%  '

$ws.Range("D44").Value = "'16.01"
$ws.Range("E44").Value = '  +3.95%  '

$ws.Range("D45").Value = "'2.582"
$ws.Range("E45").Value = '  +5.01%  '

$ws.Range("D46").Value = "'0.7126"
$ws.Range("E46").Value = '  -1.31%  '

$ws.Range("D47").Value = "'4.209"
$ws.Range("E47").Value = '  +1.56%  '

$ws.Range("E48").Value = '  +1.08%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = "'139.15"
$ws.Range("E49").Value = '  -1.05%  '

$ws.Range("B50").Value = 'Flow'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D50").Value = "'1.312"
$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("D51").Value = "'0.07958"
$ws.Range("E51").Value = '  -0.87%  '
